# Change "kubectl create -f ..." to "kubectl apply -f ..." on the two
# slides that show it ("Volume - Secret" and "Volume - ConfigMap"), and
# re-split the existing runs along the same word boundaries a live
# PowerPoint proofing/edit pass would leave behind (one run per "word":
# "$ ", "kubectl", the verb/flags, the two halves of the filename, etc).
#
# Helper: force the host to materialise an isolated <a:r> covering
# [start, start+len) by re-asserting a formatting attribute that is
# already true for that span. Re-asserting a same-valued property is
# enough to make the run split there without altering the formatting.
function Split-Run {
    param($TextRange, [int]$Start, [int]$Len)
    $sub = $TextRange.Characters($Start, $Len)
    $sub.Font.Bold = $sub.Font.Bold
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 11 : "Volume - Secret" -> CustomShape 3 -> my-secret.yaml
# ---------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$sh11 = $s11.Shapes.Item(4)
$tr11 = $sh11.TextFrame.TextRange

# Paragraph 1: "$ kubectl create -f my-secret.yaml" -> "... apply -f ..."
# Replace "create " (7 chars) with "apply " (6 chars) in one shot so the
# replacement text itself becomes the final "apply " run.
$tr11.Characters(11, 7).Text = "apply "

Split-Run $tr11 1 2        # "$ "
Split-Run $tr11 3 7        # "kubectl"
Split-Run $tr11 10 1       # " "
Split-Run $tr11 17 3       # "-f "
Split-Run $tr11 20 3       # "my-"
Split-Run $tr11 23 11      # "secret.yaml"

# Paragraph 2: "$ kubectl create secret generic mysecret --from-file=secret.key"
# (text unchanged, only re-split into finer runs)
Split-Run $tr11 35 2       # "$ "
Split-Run $tr11 37 7       # "kubectl"
Split-Run $tr11 44 23      # " create secret generic "
Split-Run $tr11 67 8       # "mysecret"
Split-Run $tr11 75 13      # " --from-file="
Split-Run $tr11 88 10      # "secret.key"

# ---------------------------------------------------------------
# Slide 9 : "Volume - ConfigMap" -> CustomShape 3 -> my-config.yaml
# ---------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(5)
$tr9 = $sh9.TextFrame.TextRange

# Paragraph 1: "$ kubectl create -f my-config.yaml" -> "... apply -f ..."
$tr9.Characters(11, 7).Text = "apply "

Split-Run $tr9 1 2         # "$ "
Split-Run $tr9 3 7         # "kubectl"
Split-Run $tr9 17 3        # "-f "
Split-Run $tr9 20 3        # "my-"
Split-Run $tr9 23 11       # "config.yaml"

# Paragraph 2: "$ kubectl create configmap --from-file=example.property.file"
# (text unchanged, only re-split into finer runs)
Split-Run $tr9 35 2        # "$ "
Split-Run $tr9 37 7        # "kubectl"
Split-Run $tr9 44 8        # " create "
Split-Run $tr9 52 9        # "configmap"
Split-Run $tr9 61 13       # " --from-file="
Split-Run $tr9 74 21       # "example.property.file"
